# "Dataframe ST.xlsx" — add a new day ("02-nov") of lookup data.
#
# Sheet3 (the lookup table) gets a fresh column of raw per-product values in
# A20:B36 (column B updated in place — same products, new readings for the
# new day). Sheet1's CB/CC columns are live VLOOKUP formulas against that
# table, so they recalc automatically. Sheet1 also keeps a frozen,
# values-only history column per day (CD=30-oct, CE=31-oct, CF=01-nov); we
# add the next one, CG=02-nov, as a values-only snapshot of the freshly
# recalculated CB column.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# 1) New day's raw values in Sheet3!B20:B36 (A20:A36 product keys unchanged).
$ws3.Range("B20").Value2 = 13.456241942693559
$ws3.Range("B21").Value2 = 1.7525390803943737
$ws3.Range("B22").Value2 = 8.6477451095384463
$ws3.Range("B23").Value2 = 5.9133386724631194
$ws3.Range("B24").Value2 = 10.508365400476409
$ws3.Range("B25").Value2 = 0
$ws3.Range("B26").Value2 = 3.9134548346044098
$ws3.Range("B27").Value2 = 7.1357289184805728
$ws3.Range("B28").Value2 = 7.0164811917716952
$ws3.Range("B29").Value2 = 5.9617133897615
$ws3.Range("B30").Value2 = 0
$ws3.Range("B31").Value2 = 0.12903874681946309
$ws3.Range("B32").Value2 = 10.192396789047077
$ws3.Range("B33").Value2 = 3.8728842971732478
$ws3.Range("B34").Value2 = 3.156161120655943
$ws3.Range("B35").Value2 = 9.6587287816936396
$ws3.Range("B36").Value2 = 43.274218130082168

# 2) New header for the day on Sheet1 (column CG), same style as the prior
#    day's frozen column (CF).
$ws1.Range("CG1").Value = "02-nov"
$ws1.Range("CG1").NumberFormat = $ws1.Range("CF1").NumberFormat

# 3) Freeze today's (now recalculated) CB values into CG2:CG18 as plain
#    numbers, matching how CD/CE/CF were produced on previous days.
for ($r = 2; $r -le 18; $r++) {
    $src = $ws1.Cells.Item($r, 80)   # column CB
    $dst = $ws1.Cells.Item($r, 85)   # column CG
    $dst.Value2 = $src.Value2
    $dst.NumberFormat = $ws1.Cells.Item($r, 84).NumberFormat   # match column CF
}

# 4) Leave the selection on the newly added header cell, like the source
#    workbook.
$ws1.Range("CG1").Select()
